# Daily update at 8 AM UTC: append the next day's win counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last data row (currently row 8) carries a distinct "last row" date
# style. Move that style down to the new last row (9) first, then restore
# row 8 to the regular interior-row date style (matching row 7).
$ws.Range("A9").NumberFormat = $ws.Range("A8").NumberFormat
$ws.Range("A8").NumberFormat = $ws.Range("A7").NumberFormat

# Append the new day's data.
$ws.Range("A9").Value = 45958
$ws.Range("B9").Value = 17
$ws.Range("C9").Value = 22
$ws.Range("D9").Value = 19
